$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Tnfsf13"
$ws.Cells.Item(2,3).Value = "Tnfrsf1a"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.637903
$ws.Cells.Item(2,8).Value = 1.913709
$ws.Cells.Item(2,9).Value = 0.1229013127714845
$ws.Cells.Item(2,10).Value = 0.1229013127714844
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 34.53682066666666
$ws.Cells.Item(2,14).Value = 103.610462
$ws.Cells.Item(2,15).Value = 0.2803141013583512
$ws.Cells.Item(2,16).Value = 0.2803141013583513
$ws.Cells.Item(2,17).Value = 22.03114151372866
$ws.Cells.Item(2,18).Value = 198.280273623558
$ws.Cells.Item(2,19).Value = 0.03445097104530032
$ws.Cells.Item(2,20).Value = 0.03445097104530032

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Tnfsf13"
$ws.Cells.Item(3,3).Value = "Tnfrsf1a"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.637903
$ws.Cells.Item(3,8).Value = 1.913709
$ws.Cells.Item(3,9).Value = 0.1229013127714845
$ws.Cells.Item(3,10).Value = 0.1229013127714844
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 47.750315
$ws.Cells.Item(3,14).Value = 143.250945
$ws.Cells.Item(3,15).Value = 0.3875598963781245
$ws.Cells.Item(3,16).Value = 0.3875598963781245
$ws.Cells.Item(3,17).Value = 30.460069189445
$ws.Cells.Item(3,18).Value = 274.140622705005
$ws.Cells.Item(3,19).Value = 0.04763162004245199
$ws.Cells.Item(3,20).Value = 0.04763162004245198

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Tnfsf13"
$ws.Cells.Item(4,3).Value = "Tnfrsf1a"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.637903
$ws.Cells.Item(4,8).Value = 1.913709
$ws.Cells.Item(4,9).Value = 0.1229013127714845
$ws.Cells.Item(4,10).Value = 0.1229013127714844
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 28.69151433333333
$ws.Cells.Item(4,14).Value = 86.074543
$ws.Cells.Item(4,15).Value = 0.232871350104353
$ws.Cells.Item(4,16).Value = 0.232871350104353
$ws.Cells.Item(4,17).Value = 18.30240306777633
$ws.Cells.Item(4,18).Value = 164.721627609987
$ws.Cells.Item(4,19).Value = 0.02862019463469295
$ws.Cells.Item(4,20).Value = 0.02862019463469295

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Tnfsf13"
$ws.Cells.Item(5,3).Value = "Tnfrsf1a"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.637903
$ws.Cells.Item(5,8).Value = 1.913709
$ws.Cells.Item(5,9).Value = 0.1229013127714845
$ws.Cells.Item(5,10).Value = 0.1229013127714844
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 12.228925
$ws.Cells.Item(5,14).Value = 36.686775
$ws.Cells.Item(5,15).Value = 0.09925465215917123
$ws.Cells.Item(5,16).Value = 0.09925465215917123
$ws.Cells.Item(5,17).Value = 7.800867944275002
$ws.Cells.Item(5,18).Value = 70.20781149847501
$ws.Cells.Item(5,19).Value = 0.0121985270490392
$ws.Cells.Item(5,20).Value = 0.0121985270490392

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Tnfsf13"
$ws.Cells.Item(6,3).Value = "Tnfrsf1a"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.170281
$ws.Cells.Item(6,8).Value = 0.510843
$ws.Cells.Item(6,9).Value = 0.03280711713229307
$ws.Cells.Item(6,10).Value = 0.03280711713229307
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 34.53682066666666
$ws.Cells.Item(6,14).Value = 103.610462
$ws.Cells.Item(6,15).Value = 0.2803141013583512
$ws.Cells.Item(6,16).Value = 0.2803141013583513
$ws.Cells.Item(6,17).Value = 5.880964359940667
$ws.Cells.Item(6,18).Value = 52.92867923946601
$ws.Cells.Item(6,19).Value = 0.0091962975570969
$ws.Cells.Item(6,20).Value = 0.009196297557096902

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Tnfsf13"
$ws.Cells.Item(7,3).Value = "Tnfrsf1a"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.170281
$ws.Cells.Item(7,8).Value = 0.510843
$ws.Cells.Item(7,9).Value = 0.03280711713229307
$ws.Cells.Item(7,10).Value = 0.03280711713229307
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 47.750315
$ws.Cells.Item(7,14).Value = 143.250945
$ws.Cells.Item(7,15).Value = 0.3875598963781245
$ws.Cells.Item(7,16).Value = 0.3875598963781245
$ws.Cells.Item(7,17).Value = 8.130971388515
$ws.Cells.Item(7,18).Value = 73.17874249663501
$ws.Cells.Item(7,19).Value = 0.01271472291625649
$ws.Cells.Item(7,20).Value = 0.01271472291625649

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Tnfsf13"
$ws.Cells.Item(8,3).Value = "Tnfrsf1a"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.170281
$ws.Cells.Item(8,8).Value = 0.510843
$ws.Cells.Item(8,9).Value = 0.03280711713229307
$ws.Cells.Item(8,10).Value = 0.03280711713229307
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 28.69151433333333
$ws.Cells.Item(8,14).Value = 86.074543
$ws.Cells.Item(8,15).Value = 0.232871350104353
$ws.Cells.Item(8,16).Value = 0.232871350104353
$ws.Cells.Item(8,17).Value = 4.885619752194334
$ws.Cells.Item(8,18).Value = 43.97057776974901
$ws.Cells.Item(8,19).Value = 0.007639837659628738
$ws.Cells.Item(8,20).Value = 0.007639837659628738

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Tnfsf13"
$ws.Cells.Item(9,3).Value = "Tnfrsf1a"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.170281
$ws.Cells.Item(9,8).Value = 0.510843
$ws.Cells.Item(9,9).Value = 0.03280711713229307
$ws.Cells.Item(9,10).Value = 0.03280711713229307
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 12.228925
$ws.Cells.Item(9,14).Value = 36.686775
$ws.Cells.Item(9,15).Value = 0.09925465215917123
$ws.Cells.Item(9,16).Value = 0.09925465215917123
$ws.Cells.Item(9,17).Value = 2.082353577925001
$ws.Cells.Item(9,18).Value = 18.741182201325
$ws.Cells.Item(9,19).Value = 0.003256258999310936
$ws.Cells.Item(9,20).Value = 0.003256258999310936

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Tnfsf13"
$ws.Cells.Item(10,3).Value = "Tnfrsf1a"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 3.407124666666667
$ws.Cells.Item(10,8).Value = 10.221374
$ws.Cells.Item(10,9).Value = 0.6564322386153376
$ws.Cells.Item(10,10).Value = 0.6564322386153377
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 34.53682066666666
$ws.Cells.Item(10,14).Value = 103.610462
$ws.Cells.Item(10,15).Value = 0.2803141013583512
$ws.Cells.Item(10,16).Value = 0.2803141013583513
$ws.Cells.Item(10,17).Value = 117.6712536016431
$ws.Cells.Item(10,18).Value = 1059.041282414788
$ws.Cells.Item(10,19).Value = 0.1840072130701091
$ws.Cells.Item(10,20).Value = 0.1840072130701092

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Tnfsf13"
$ws.Cells.Item(11,3).Value = "Tnfrsf1a"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 3.407124666666667
$ws.Cells.Item(11,8).Value = 10.221374
$ws.Cells.Item(11,9).Value = 0.6564322386153376
$ws.Cells.Item(11,10).Value = 0.6564322386153377
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 47.750315
$ws.Cells.Item(11,14).Value = 143.250945
$ws.Cells.Item(11,15).Value = 0.3875598963781245
$ws.Cells.Item(11,16).Value = 0.3875598963781245
$ws.Cells.Item(11,17).Value = 162.6912760776033
$ws.Cells.Item(11,18).Value = 1464.22148469843
$ws.Cells.Item(11,19).Value = 0.2544068103770205
$ws.Cells.Item(11,20).Value = 0.2544068103770206

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Tnfsf13"
$ws.Cells.Item(12,3).Value = "Tnfrsf1a"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 3.407124666666667
$ws.Cells.Item(12,8).Value = 10.221374
$ws.Cells.Item(12,9).Value = 0.6564322386153376
$ws.Cells.Item(12,10).Value = 0.6564322386153377
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 28.69151433333333
$ws.Cells.Item(12,14).Value = 86.074543
$ws.Cells.Item(12,15).Value = 0.232871350104353
$ws.Cells.Item(12,16).Value = 0.232871350104353
$ws.Cells.Item(12,17).Value = 97.75556620912023
$ws.Cells.Item(12,18).Value = 879.8000958820821
$ws.Cells.Item(12,19).Value = 0.1528642616583765
$ws.Cells.Item(12,20).Value = 0.1528642616583765

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Tnfsf13"
$ws.Cells.Item(13,3).Value = "Tnfrsf1a"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 3.407124666666667
$ws.Cells.Item(13,8).Value = 10.221374
$ws.Cells.Item(13,9).Value = 0.6564322386153376
$ws.Cells.Item(13,10).Value = 0.6564322386153377
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 12.228925
$ws.Cells.Item(13,14).Value = 36.686775
$ws.Cells.Item(13,15).Value = 0.09925465215917123
$ws.Cells.Item(13,16).Value = 0.09925465215917123
$ws.Cells.Item(13,17).Value = 41.66547201431668
$ws.Cells.Item(13,18).Value = 374.9892481288501
$ws.Cells.Item(13,19).Value = 0.06515395350983143
$ws.Cells.Item(13,20).Value = 0.06515395350983144

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Tnfsf13"
$ws.Cells.Item(14,3).Value = "Tnfrsf1a"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 0.6666666666666666
$ws.Cells.Item(14,7).Value = 0.975059
$ws.Cells.Item(14,8).Value = 2.925177
$ws.Cells.Item(14,9).Value = 0.1878593314808848
$ws.Cells.Item(14,10).Value = 0.1878593314808848
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 34.53682066666666
$ws.Cells.Item(14,14).Value = 103.610462
$ws.Cells.Item(14,15).Value = 0.2803141013583512
$ws.Cells.Item(14,16).Value = 0.2803141013583513
$ws.Cells.Item(14,17).Value = 33.67543782241933
$ws.Cells.Item(14,18).Value = 303.078940401774
$ws.Cells.Item(14,19).Value = 0.05265961968584484
$ws.Cells.Item(14,20).Value = 0.05265961968584486

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Tnfsf13"
$ws.Cells.Item(15,3).Value = "Tnfrsf1a"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 0.6666666666666666
$ws.Cells.Item(15,7).Value = 0.975059
$ws.Cells.Item(15,8).Value = 2.925177
$ws.Cells.Item(15,9).Value = 0.1878593314808848
$ws.Cells.Item(15,10).Value = 0.1878593314808848
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 47.750315
$ws.Cells.Item(15,14).Value = 143.250945
$ws.Cells.Item(15,15).Value = 0.3875598963781245
$ws.Cells.Item(15,16).Value = 0.3875598963781245
$ws.Cells.Item(15,17).Value = 46.559374393585
$ws.Cells.Item(15,18).Value = 419.034369542265
$ws.Cells.Item(15,19).Value = 0.07280674304239546
$ws.Cells.Item(15,20).Value = 0.07280674304239546

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Tnfsf13"
$ws.Cells.Item(16,3).Value = "Tnfrsf1a"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 0.6666666666666666
$ws.Cells.Item(16,7).Value = 0.975059
$ws.Cells.Item(16,8).Value = 2.925177
$ws.Cells.Item(16,9).Value = 0.1878593314808848
$ws.Cells.Item(16,10).Value = 0.1878593314808848
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 28.69151433333333
$ws.Cells.Item(16,14).Value = 86.074543
$ws.Cells.Item(16,15).Value = 0.232871350104353
$ws.Cells.Item(16,16).Value = 0.232871350104353
$ws.Cells.Item(16,17).Value = 27.97591927434567
$ws.Cells.Item(16,18).Value = 251.783273469111
$ws.Cells.Item(16,19).Value = 0.04374705615165483
$ws.Cells.Item(16,20).Value = 0.04374705615165484

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Tnfsf13"
$ws.Cells.Item(17,3).Value = "Tnfrsf1a"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 2
$ws.Cells.Item(17,6).Value = 0.6666666666666666
$ws.Cells.Item(17,7).Value = 0.975059
$ws.Cells.Item(17,8).Value = 2.925177
$ws.Cells.Item(17,9).Value = 0.1878593314808848
$ws.Cells.Item(17,10).Value = 0.1878593314808848
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 12.228925
$ws.Cells.Item(17,14).Value = 36.686775
$ws.Cells.Item(17,15).Value = 0.09925465215917123
$ws.Cells.Item(17,16).Value = 0.09925465215917123
$ws.Cells.Item(17,17).Value = 11.923923381575
$ws.Cells.Item(17,18).Value = 107.315310434175
$ws.Cells.Item(17,19).Value = 0.01864591260098966
$ws.Cells.Item(17,20).Value = 0.01864591260098967
